# Repull data, push all data, mean calculation
# Updates the dSF (column F) values for the rows whose underlying data
# was re-pulled, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    11 = -1
    12 = 0
    16 = 3
    21 = 2
    24 = 3
    26 = 3
    27 = -4
    34 = 1
    36 = 1
    39 = 1
    45 = -4
    49 = -7
    50 = 3
    58 = 3
    61 = 1
    62 = -3
    64 = -4
    65 = 0
    67 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
